$d = $word.ActiveDocument

# 1. Tighten the right paragraph border spacing (w:space 16 -> 5) on every
#    paragraph that currently uses the old value (the four bordered
#    "boxed" paragraphs in the template).
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Borders.DistanceFromRight -eq 16) {
        $p.Borders.DistanceFromRight = 5
    }
}

# 2. Relocate the stray "_GoBack" bookmark (an artifact of the last cursor
#    position) from the middle of the letterhead down to the very end of
#    the document, where Word normally leaves it after the final edit.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$lastParagraph = $d.Paragraphs.Last
$lastParagraph.Range.Bookmarks.Add("_GoBack")
